# création vue initialisation projet
# Update column A (date values) from 2015xxxx -> 2017xxxx for rows 3..63
# and update specific column E values per the target revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: shift the "year" portion of the date-like integer by +20000
# (20150926 -> 20170926, etc.) for every data row (3 through 63).
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2
    if ($null -ne $oldVal) {
        $cell.Value2 = [double]$oldVal + 20000
    }
}

# Column E: explicit updated values as per the diff.
$eUpdates = @{
    4  = 10
    6  = 17
    7  = 19
    8  = 8
    9  = 6
    10 = 9
    11 = 17
    12 = 19
    13 = 10
    14 = 12
    15 = 7
    16 = 19
    17 = 5
    18 = 17
    19 = 5
    20 = 13
    21 = 16
    22 = 17
    23 = 10
    25 = 20
    26 = 12
    27 = 19
    28 = 14
    29 = 10
    30 = 11
    31 = 7
    32 = 15
    33 = 11
    34 = 8
    35 = 18
    36 = 10
    37 = 10
    38 = 20
    39 = 10
    40 = 8
    41 = 16
    42 = 18
    43 = 12
    44 = 7
    45 = 19
    46 = 5
    47 = 14
    48 = 9
    49 = 15
    50 = 18
    51 = 16
    52 = 14
    53 = 6
    54 = 17
    55 = 13
    56 = 8
    57 = 17
    58 = 6
    59 = 18
    60 = 12
    61 = 16
    62 = 12
    63 = 6
}

foreach ($row in $eUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value2 = $eUpdates[$row]
}
